$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Set Runmode (column D) to "Y" for all test case rows (2 through 25)
# i.e. run all the test cases
$ws.Range("D2:D25").Value = "Y"

# Reflect the resulting selection left behind by this edit
$ws.Range("D2:D25").Select()
